{"js": "// Replace the three-digit-divided-by-one-digit problems in the table\n// with the new set of problems, per the commit's regenerated values.\nconst replacements = [\n  [\"647\u00f74=\", \"255\u00f75=\"],\n  [\"635\u00f74=\", \"848\u00f75=\"],\n  [\"358\u00f79=\", \"702\u00f76=\"],\n  [\"565\u00f74=\", \"916\u00f78=\"],\n  [\"578\u00f73=\", \"422\u00f73=\"],\n  [\"227\u00f76=\", \"242\u00f77=\"],\n  [\"857\u00f73=\", \"400\u00f73=\"],\n  [\"539\u00f76=\", \"155\u00f77=\"],\n  [\"384\u00f79=\", \"971\u00f74=\"],\n  [\"432\u00f74=\", \"135\u00f77=\"],\n  [\"751\u00f77=\", \"928\u00f76=\"],\n  [\"677\u00f74=\", \"101\u00f76=\"],\n  [\"152\u00f74=\", \"656\u00f74=\"],\n  [\"695\u00f73=\", \"157\u00f78=\"],\n  [\"137\u00f78=\", \"759\u00f76=\"],\n  [\"591\u00f75=\", \"877\u00f77=\"],\n  [\"127\u00f72=\", \"153\u00f72=\"],\n  [\"625\u00f79=\", \"937\u00f79=\"],\n  [\"146\u00f76=\", \"150\u00f76=\"],\n  [\"834\u00f74=\", \"962\u00f77=\"],\n  [\"404\u00f72=\", \"702\u00f74=\"],\n  [\"657\u00f77=\", \"578\u00f75=\"],\n  [\"421\u00f78=\", \"270\u00f77=\"],\n  [\"623\u00f74=\", \"698\u00f73=\"],\n  [\"178\u00f77=\", \"299\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-divided-by-one-digit problems in the table\n# with the new set of problems, per the commit's regenerated values.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @(\"647\u00f74=\", \"255\u00f75=\"),\n    @(\"635\u00f74=\", \"848\u00f75=\"),\n    @(\"358\u00f79=\", \"702\u00f76=\"),\n    @(\"565\u00f74=\", \"916\u00f78=\"),\n    @(\"578\u00f73=\", \"422\u00f73=\"),\n    @(\"227\u00f76=\", \"242\u00f77=\"),\n    @(\"857\u00f73=\", \"400\u00f73=\"),\n    @(\"539\u00f76=\", \"155\u00f77=\"),\n    @(\"384\u00f79=\", \"971\u00f74=\"),\n    @(\"432\u00f74=\", \"135\u00f77=\"),\n    @(\"751\u00f77=\", \"928\u00f76=\"),\n    @(\"677\u00f74=\", \"101\u00f76=\"),\n    @(\"152\u00f74=\", \"656\u00f74=\"),\n    @(\"695\u00f73=\", \"157\u00f78=\"),\n    @(\"137\u00f78=\", \"759\u00f76=\"),\n    @(\"591\u00f75=\", \"877\u00f77=\"),\n    @(\"127\u00f72=\", \"153\u00f72=\"),\n    @(\"625\u00f79=\", \"937\u00f79=\"),\n    @(\"146\u00f76=\", \"150\u00f76=\"),\n    @(\"834\u00f74=\", \"962\u00f77=\"),\n    @(\"404\u00f72=\", \"702\u00f74=\"),\n    @(\"657\u00f77=\", \"578\u00f75=\"),\n    @(\"421\u00f78=\", \"270\u00f77=\"),\n    @(\"623\u00f74=\", \"698\u00f73=\"),\n    @(\"178\u00f77=\", \"299\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll)\n}\n"}
